$wb = $excel.ActiveWorkbook

# Worksheets
$wsApplicants = $wb.Worksheets.Item("Applicants")
$wsJobOffers  = $wb.Worksheets.Item("JobOffers")
$wsSkills     = $wb.Worksheets.Item("Skills")

# --- JobOffers sheet: add Level/Skills columns (E:H) ---
$wsJobOffers.Range("E1").Value = "Level"
$wsJobOffers.Range("F1").Value = "Skills"
$wsJobOffers.Range("G1").Value = "Skills"
$wsJobOffers.Range("H1").Value = "Skills"

$wsJobOffers.Range("E2").Value = "Junior"
$wsJobOffers.Range("F2").Value = "Java fundamentals"
$wsJobOffers.Range("G2").Value = "Java Spring"

$wsJobOffers.Range("E3").Value = "Junior"
$wsJobOffers.Range("F3").Value = "Graphics"
$wsJobOffers.Range("G3").Value = "Multimedia"

$wsJobOffers.Range("E4").Value = "Mid"
$wsJobOffers.Range("F4").Value = "Hardware"
$wsJobOffers.Range("G4").Value = "Devops"
$wsJobOffers.Range("H4").Value = "Databases"

$wsJobOffers.Range("E5").Value = "Senior"
$wsJobOffers.Range("F5").Value = "Java Spring"

$wsJobOffers.Range("E6").Value = "Senior"
$wsJobOffers.Range("F6").Value = "Scrum"

$wsJobOffers.Range("E7").Value = "Senior"
$wsJobOffers.Range("F7").Value = "Problem-solving"

$wsJobOffers.Range("E8").Value = "Junior"
$wsJobOffers.Range("F8").Value = "Java Spring"
$wsJobOffers.Range("G8").Value = "Devops"

$wsJobOffers.Range("E9").Value = "Mid"
$wsJobOffers.Range("F9").Value = "Operating Systems"
$wsJobOffers.Range("G9").Value = "Angular"
$wsJobOffers.Range("H9").Value = "C#"

$wsJobOffers.Range("E10").Value = "Mid"
$wsJobOffers.Range("F10").Value = "Maven"
$wsJobOffers.Range("G10").Value = "Spreadsheets"

$wsJobOffers.Range("E11").Value = "Junior"
$wsJobOffers.Range("F11").Value = "Maven"
$wsJobOffers.Range("G11").Value = "Problem-solving"

$wsJobOffers.Range("E12").Value = "Junior"
$wsJobOffers.Range("F12").Value = "Angular"
$wsJobOffers.Range("G12").Value = "C#"

$wsJobOffers.Range("E13").Value = "Senior"
$wsJobOffers.Range("F13").Value = "Operating Systems"

$wsJobOffers.Range("E14").Value = "Mid"
$wsJobOffers.Range("F14").Value = "Devops"
$wsJobOffers.Range("G14").Value = "Problem-solving"
$wsJobOffers.Range("H14").Value = "Databases"

# Header style to match A1:D1 (bold) for new header cells
$wsJobOffers.Range("E1:H1").Font.Bold = $true

# --- Selections / views ---
$wsApplicants.Activate()
$wsApplicants.Range("G2:G10").Select()

$wsJobOffers.Activate()
$wsJobOffers.Range("E7").Select()

$wsSkills.Activate()
$wsSkills.Range("A21").Select()

# Make JobOffers the active tab as final state
$wsJobOffers.Activate()
